$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 additions
$ws.Range("B2").Value = "rekisela@uw.edu"
$ws.Range("C2").Value = "Jim Kisela"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = "Ace Ventura"

# Row 5 additions
$ws.Range("B5").Value = "rekisela@uw.edu"
$ws.Range("C5").Value = "Napoleon Bonaparte"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "No"
$ws.Range("H5").Value = "Ace Ventura"
$ws.Range("I5").Value = "Blazing Saddles"
$ws.Range("J5").Value = "Contact"
$ws.Range("K5").Value = "Editor"
$ws.Range("L5").Value = "Special Effects"
$ws.Range("M5").Value = "Camera Operator"
$ws.Range("N5").Value = "Role"

# Row 9 additions
$ws.Range("B9").Value = "rekisela@uw.edu"
$ws.Range("C9").Value = "Ann Kisela"
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "No"
$ws.Range("H9").Value = "Ace Ventura"
$ws.Range("I9").Value = "Contact"
$ws.Range("J9").Value = "Blazing Saddles"
$ws.Range("K9").Value = "Camera Operator"
$ws.Range("L9").Value = "Assistant Director"
$ws.Range("M9").Value = "Production Assistant"
$ws.Range("N9").Value = "Role"

# New Row 16
$ws.Range("A16").Value = "Production Assistant"
$ws.Range("B16").Value = "rekisela@uw.edu"
$ws.Range("C16").Value = "Fake Name 3"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = "No"
$ws.Range("H16").Value = "Ace Ventura"
$ws.Range("I16").Value = "Blazing Saddles"
$ws.Range("J16").Value = "Contact"
$ws.Range("K16").Value = "Production Assistant"
$ws.Range("L16").Value = "Assistant Camera Operator"
$ws.Range("M16").Value = "Art Department"
$ws.Range("N16").Value = "Role"
